$d = $word.ActiveDocument

# Unicode helpers (curly apostrophe / em dash) so the script stays plain ASCII.
$rsquo = [char]0x2019
$mdash = [char]0x2014

# 1) Section heading: "THREE PRINCIPLES" -> "TWO PRINCIPLES"
$d.Content.Find.Execute("THREE PRINCIPLES", $true, $false, $false, $false, $false, $true, 1, $false, "TWO PRINCIPLES", 2) | Out-Null

# 2) Reorder the five ASSEMBLY bullet paragraphs (the five paragraphs right after
#    the "ASSEMBLY" heading) in place.
#    Before: One package..., When it's empty..., Wash hands..., Gloves required..., Hair tied back...
#    After : Hair tied back..., Wash hands..., Gloves required..., One package..., When it's empty...
#    Each bullet paragraph is "<teal bullet run>• <gray text run>". Rewriting just the
#    text of the gray run (leaving both runs' formatting untouched) achieves the reorder
#    without disturbing any paragraph/run properties.
$bulletTexts = @(
    "Hair tied back or covered" + $mdash + " hair and beard nets are best.",
    "Wash hands with soap and water" + $mdash + " hand sanitizer doesn" + $rsquo + "t kill some food-borne illness-causing organisms.",
    "Gloves required " + $mdash + " change them after touching phone, face, other surfaces.",
    "One package each of meat and cheese out at a time per table.",
    "When it" + $rsquo + "s empty, return to fridge for the next one. This is where runners come in handy."
)

# Locate the "ASSEMBLY" heading paragraph; the five bullets immediately follow it.
$assemblyIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.Trim() -eq "ASSEMBLY") {
        $assemblyIndex = $i
        break
    }
}

$startIndex = $assemblyIndex + 1
for ($i = 0; $i -lt $bulletTexts.Length; $i++) {
    $p = $d.Paragraphs.Item($startIndex + $i)
    $r = $p.Range
    # Skip the leading two-character bullet run ("• ") and the trailing paragraph mark.
    $textRange = $d.Range($r.Start + 2, $r.End - 1)
    $textRange.Text = $bulletTexts[$i]
}

# 3) Section heading: "BUILDING THE SANDWICH" -> "MAKING THE SANDWICHES"
$d.Content.Find.Execute("BUILDING THE SANDWICH", $true, $false, $false, $false, $false, $true, 1, $false, "MAKING THE SANDWICHES", 2) | Out-Null

# 4) Header title: "Food Safety Guidelines" -> "Group Event Food Safety Guidelines"
$d.Sections.Item(1).Headers.Item(1).Range.Find.Execute("Food Safety Guidelines", $true, $false, $false, $false, $false, $true, 1, $false, "Group Event Food Safety Guidelines", 2) | Out-Null
